$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 50 and 51 swap their Coin/Link/Price/Volume content
# (Algorand <-> RocketPoolETH), plus the Price/Volume values change.
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.670.22"
$ws.Range("E50").Value = "  +7.05%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.187"
$ws.Range("E51").Value = "  +0.56%  "

$ws.Range("D2").Value = "46.403.17"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.430.90"
$ws.Range("E3").Value = "  +6.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.37"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.43"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.32"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.17"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").Value = "2.801.62"
$ws.Range("E14").Value = "  +7.02%  "
$ws.Range("D15").Value = "2.421.58"
$ws.Range("E15").Value = "  +6.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("E16").Value = "  +6.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.15"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").Value = "46.256.69"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.24"
$ws.Range("E21").Value = "  +6.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.64"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.56"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +5.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.41"
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  +14.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.41"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("E33").Value = "  +3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.34"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0774"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.02"
$ws.Range("E36").Value = "  +18.66%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.18"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.93"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0303"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("D43").Value = "1.985.31"
$ws.Range("E43").Value = "  +11.25%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.78"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.84"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.48"
$ws.Range("E47").Value = "  +32.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.66"
$ws.Range("E48").Value = "  +10.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.30"
$ws.Range("E49").Value = "  +6.71%  "
